$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Qty executed upto date" (C) column values (numeric)
$ws.Range("C8").Value = 29
$ws.Range("C9").Value = 60
$ws.Range("C10").Value = 50
$ws.Range("C11").Value = 95
$ws.Range("C12").Value = 87
$ws.Range("C13").Value = 90
$ws.Range("C14").Value = 38
$ws.Range("C15").Value = 84
$ws.Range("C16").Value = 55
$ws.Range("C17").Value = 37

# Update the "Upto date Amount" (G) column values - these are stored as
# text (e.g. "15360.00"), so force a Text number format on each target
# cell before assigning the value, to avoid Excel auto-converting the
# numeric-looking string into a real number.
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "15360.00"

$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "23600.00"

$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "62890.00"

$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "12240.00"

$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "874.00"

# Update the Grand Total rows to reflect the new sum
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = "114964.00"

$ws.Range("H19").NumberFormat = "@"
$ws.Range("H19").Value = "114964.00"

$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = "114964.00"

$ws.Range("H21").NumberFormat = "@"
$ws.Range("H21").Value = "114964.00"
